$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIN")

# Row 128: add combined-budget ref id
$ws.Range("E128").Value = 46000000000009

# Row 157: clear the user-role ref id (data moved down to row 158)
$ws.Range("D157").Value = ""

# Row 158: set user-role ref id + combined-budget ref id
#  (set E before D so the dependent CONCATENATE formula in H158 picks up
#   the fresh E value on the same recalc pass)
$ws.Range("E158").Value = 46000000000033
$ws.Range("D158").Value = 95000000000023

# Row 235: add combined-budget ref id
$ws.Range("E235").Value = 46000000000033

# Row 281: set user-role ref id + combined-budget ref id (E before D)
$ws.Range("E281").Value = 46000000000009
$ws.Range("D281").Value = 95000000000018

# Row 545: add combined-budget ref id
$ws.Range("E545").Value = 46000000000033

# Row 562: set user-role ref id
$ws.Range("D562").Value = 95000000000015
